$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32 (shifts existing rows 32-63 down to 33-64),
# mirroring the weekly price data being pushed back in time by one entry.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new weekly observation.
$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44589
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 100112030
$ws.Range("G32").Value = "Poroto granado"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 560
$ws.Range("K32").Value = 29000
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = 29500
$ws.Range("N32").Value = "`$/malla 25 kilos"
$ws.Range("O32").Value = "Provincia del Elquí"
$ws.Range("P32").Value = 1180
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
